$d = $word.ActiveDocument

# Find the "Delete Products" list item (Seller > Catalogue > Delete Products)
# so a new sibling item ("List Products") can be inserted right after it,
# before "Inventory".
$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Delete Products") {
        $targetIdx = $idx
        break
    }
    $idx = $idx + 1
}

if ($targetIdx -lt 0) {
    throw "Could not locate the 'Delete Products' paragraph."
}

$target = $d.Paragraphs.Item($targetIdx + 1)

# Splitting right after "Delete Products" creates a new paragraph that
# inherits the same pPr/rPr (ListParagraph style, ilvl=2/numId=1, Times New
# Roman / 333333 / white-shaded run) as the paragraph it was split from.
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIdx + 2)
$insertionPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertionPoint.InsertAfter("List Products")
